$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = [double]"10.587229"
$ws.Range("H2").Value2 = [double]"31.761687"
$ws.Range("I2").Value2 = [double]"0.1340590927938227"
$ws.Range("J2").Value2 = [double]"0.1340590927938227"
$ws.Range("M2").Value2 = [double]"2.429935333333333"
$ws.Range("N2").Value2 = [double]"7.289806"
$ws.Range("O2").Value2 = [double]"0.0007999299918632063"
$ws.Range("P2").Value2 = [double]"0.0007999299918632063"
$ws.Range("Q2").Value2 = [double]"25.72628182919134"
$ws.Range("R2").Value2 = [double]"231.536536462722"
$ws.Range("S2").Value2 = [double]"0.0001072378890077514"
$ws.Range("T2").Value2 = [double]"0.0001072378890077514"
$ws.Range("G3").Value2 = [double]"10.587229"
$ws.Range("H3").Value2 = [double]"31.761687"
$ws.Range("I3").Value2 = [double]"0.1340590927938227"
$ws.Range("J3").Value2 = [double]"0.1340590927938227"
$ws.Range("O3").Value2 = [double]"0.001089498967837074"
$ws.Range("P3").Value2 = [double]"0.001089498967837074"
$ws.Range("Q3").Value2 = [double]"35.03901314401867"
$ws.Range("R3").Value2 = [double]"315.351118296168"
$ws.Range("S3").Value2 = [double]"0.0001460572432280443"
$ws.Range("T3").Value2 = [double]"0.0001460572432280443"
$ws.Range("G4").Value2 = [double]"10.587229"
$ws.Range("H4").Value2 = [double]"31.761687"
$ws.Range("I4").Value2 = [double]"0.1340590927938227"
$ws.Range("J4").Value2 = [double]"0.1340590927938227"
$ws.Range("M4").Value2 = [double]"1.214141"
$ws.Range("N4").Value2 = [double]"3.642423"
$ws.Range("O4").Value2 = [double]"0.0003996928588706414"
$ws.Range("P4").Value2 = [double]"0.0003996928588706414"
$ws.Range("Q4").Value2 = [double]"12.854388805289"
$ws.Range("R4").Value2 = [double]"115.689499247601"
$ws.Range("S4").Value2 = [double]"5.358246205636758E-05"
$ws.Range("T4").Value2 = [double]"5.358246205636758E-05"
$ws.Range("G5").Value2 = [double]"10.587229"
$ws.Range("H5").Value2 = [double]"31.761687"
$ws.Range("I5").Value2 = [double]"0.1340590927938227"
$ws.Range("J5").Value2 = [double]"0.1340590927938227"
$ws.Range("M5").Value2 = [double]"3030.731364"
$ws.Range("N5").Value2 = [double]"9092.194092"
$ws.Range("O5").Value2 = [double]"0.9977108781814291"
$ws.Range("P5").Value2 = [double]"0.9977108781814292"
$ws.Range("Q5").Value2 = [double]"32087.04698815036"
$ws.Range("R5").Value2 = [double]"288783.4228933532"
$ws.Range("S5").Value2 = [double]"0.1337522151995305"
$ws.Range("T5").Value2 = [double]"0.1337522151995305"
$ws.Range("I6").Value2 = [double]"0.2241137347582675"
$ws.Range("J6").Value2 = [double]"0.2241137347582675"
$ws.Range("M6").Value2 = [double]"2.429935333333333"
$ws.Range("N6").Value2 = [double]"7.289806"
$ws.Range("O6").Value2 = [double]"0.0007999299918632063"
$ws.Range("P6").Value2 = [double]"0.0007999299918632063"
$ws.Range("Q6").Value2 = [double]"43.00799730944845"
$ws.Range("R6").Value2 = [double]"387.071975785036"
$ws.Range("S6").Value2 = [double]"0.0001792752980216137"
$ws.Range("T6").Value2 = [double]"0.0001792752980216137"
$ws.Range("I7").Value2 = [double]"0.2241137347582675"
$ws.Range("J7").Value2 = [double]"0.2241137347582675"
$ws.Range("O7").Value2 = [double]"0.001089498967837074"
$ws.Range("P7").Value2 = [double]"0.001089498967837074"
$ws.Range("S7").Value2 = [double]"0.0002441716826972442"
$ws.Range("T7").Value2 = [double]"0.0002441716826972442"
$ws.Range("I8").Value2 = [double]"0.2241137347582675"
$ws.Range("J8").Value2 = [double]"0.2241137347582675"
$ws.Range("M8").Value2 = [double]"1.214141"
$ws.Range("N8").Value2 = [double]"3.642423"
$ws.Range("O8").Value2 = [double]"0.0003996928588706414"
$ws.Range("P8").Value2 = [double]"0.0003996928588706414"
$ws.Range("Q8").Value2 = [double]"21.48936728684867"
$ws.Range("R8").Value2 = [double]"193.404305581638"
$ws.Range("S8").Value2 = [double]"8.957665935770859E-05"
$ws.Range("T8").Value2 = [double]"8.957665935770859E-05"
$ws.Range("I9").Value2 = [double]"0.2241137347582675"
$ws.Range("J9").Value2 = [double]"0.2241137347582675"
$ws.Range("M9").Value2 = [double]"3030.731364"
$ws.Range("N9").Value2 = [double]"9092.194092"
$ws.Range("O9").Value2 = [double]"0.9977108781814291"
$ws.Range("P9").Value2 = [double]"0.9977108781814292"
$ws.Range("Q9").Value2 = [double]"53641.62764355032"
$ws.Range("R9").Value2 = [double]"482774.6487919529"
$ws.Range("S9").Value2 = [double]"0.223600711118191"
$ws.Range("T9").Value2 = [double]"0.223600711118191"
$ws.Range("G10").Value2 = [double]"34.22308866666667"
$ws.Range("H10").Value2 = [double]"102.669266"
$ws.Range("I10").Value2 = [double]"0.43334438305395"
$ws.Range("J10").Value2 = [double]"0.43334438305395"
$ws.Range("M10").Value2 = [double]"2.429935333333333"
$ws.Range("N10").Value2 = [double]"7.289806"
$ws.Range("O10").Value2 = [double]"0.0007999299918632063"
$ws.Range("P10").Value2 = [double]"0.0007999299918632063"
$ws.Range("Q10").Value2 = [double]"83.15989236693289"
$ws.Range("R10").Value2 = [double]"748.4390313023961"
$ws.Range("S10").Value2 = [double]"0.0003466451688103123"
$ws.Range("T10").Value2 = [double]"0.0003466451688103123"
$ws.Range("G11").Value2 = [double]"34.22308866666667"
$ws.Range("H11").Value2 = [double]"102.669266"
$ws.Range("I11").Value2 = [double]"0.43334438305395"
$ws.Range("J11").Value2 = [double]"0.43334438305395"
$ws.Range("O11").Value2 = [double]"0.001089498967837074"
$ws.Range("P11").Value2 = [double]"0.001089498967837074"
$ws.Range("Q11").Value2 = [double]"113.2631828045138"
$ws.Range("R11").Value2 = [double]"1019.368645240624"
$ws.Range("S11").Value2 = [double]"0.000472128258055272"
$ws.Range("T11").Value2 = [double]"0.000472128258055272"
$ws.Range("G12").Value2 = [double]"34.22308866666667"
$ws.Range("H12").Value2 = [double]"102.669266"
$ws.Range("I12").Value2 = [double]"0.43334438305395"
$ws.Range("J12").Value2 = [double]"0.43334438305395"
$ws.Range("M12").Value2 = [double]"1.214141"
$ws.Range("N12").Value2 = [double]"3.642423"
$ws.Range("O12").Value2 = [double]"0.0003996928588706414"
$ws.Range("P12").Value2 = [double]"0.0003996928588706414"
$ws.Range("Q12").Value2 = [double]"41.55165509683533"
$ws.Range("R12").Value2 = [double]"373.964895871518"
$ws.Range("S12").Value2 = [double]"0.0001732046553383676"
$ws.Range("T12").Value2 = [double]"0.0001732046553383676"
$ws.Range("G13").Value2 = [double]"34.22308866666667"
$ws.Range("H13").Value2 = [double]"102.669266"
$ws.Range("I13").Value2 = [double]"0.43334438305395"
$ws.Range("J13").Value2 = [double]"0.43334438305395"
$ws.Range("M13").Value2 = [double]"3030.731364"
$ws.Range("N13").Value2 = [double]"9092.194092"
$ws.Range("O13").Value2 = [double]"0.9977108781814291"
$ws.Range("P13").Value2 = [double]"0.9977108781814292"
$ws.Range("Q13").Value2 = [double]"103720.9881950196"
$ws.Range("R13").Value2 = [double]"933488.8937551765"
$ws.Range("S13").Value2 = [double]"0.432352404971746"
$ws.Range("T13").Value2 = [double]"0.432352404971746"
$ws.Range("G14").Value2 = [double]"16.46479166666667"
$ws.Range("H14").Value2 = [double]"49.394375"
$ws.Range("I14").Value2 = [double]"0.2084827893939599"
$ws.Range("J14").Value2 = [double]"0.2084827893939599"
$ws.Range("M14").Value2 = [double]"2.429935333333333"
$ws.Range("N14").Value2 = [double]"7.289806"
$ws.Range("O14").Value2 = [double]"0.0007999299918632063"
$ws.Range("P14").Value2 = [double]"0.0007999299918632063"
$ws.Range("Q14").Value2 = [double]"40.00837902680556"
$ws.Range("R14").Value2 = [double]"360.07541124125"
$ws.Range("S14").Value2 = [double]"0.0001667716360235289"
$ws.Range("T14").Value2 = [double]"0.0001667716360235289"
$ws.Range("G15").Value2 = [double]"16.46479166666667"
$ws.Range("H15").Value2 = [double]"49.394375"
$ws.Range("I15").Value2 = [double]"0.2084827893939599"
$ws.Range("J15").Value2 = [double]"0.2084827893939599"
$ws.Range("O15").Value2 = [double]"0.001089498967837074"
$ws.Range("P15").Value2 = [double]"0.001089498967837074"
$ws.Range("Q15").Value2 = [double]"54.49112809611111"
$ws.Range("R15").Value2 = [double]"490.4201528649999"
$ws.Range("S15").Value2 = [double]"0.0002271417838565133"
$ws.Range("T15").Value2 = [double]"0.0002271417838565133"
$ws.Range("G16").Value2 = [double]"16.46479166666667"
$ws.Range("H16").Value2 = [double]"49.394375"
$ws.Range("I16").Value2 = [double]"0.2084827893939599"
$ws.Range("J16").Value2 = [double]"0.2084827893939599"
$ws.Range("M16").Value2 = [double]"1.214141"
$ws.Range("N16").Value2 = [double]"3.642423"
$ws.Range("O16").Value2 = [double]"0.0003996928588706414"
$ws.Range("P16").Value2 = [double]"0.0003996928588706414"
$ws.Range("Q16").Value2 = [double]"19.99057861895833"
$ws.Range("R16").Value2 = [double]"179.915207570625"
$ws.Range("S16").Value2 = [double]"8.332908211819767E-05"
$ws.Range("T16").Value2 = [double]"8.332908211819767E-05"
$ws.Range("G17").Value2 = [double]"16.46479166666667"
$ws.Range("H17").Value2 = [double]"49.394375"
$ws.Range("I17").Value2 = [double]"0.2084827893939599"
$ws.Range("J17").Value2 = [double]"0.2084827893939599"
$ws.Range("M17").Value2 = [double]"3030.731364"
$ws.Range("N17").Value2 = [double]"9092.194092"
$ws.Range("O17").Value2 = [double]"0.9977108781814291"
$ws.Range("P17").Value2 = [double]"0.9977108781814292"
$ws.Range("Q17").Value2 = [double]"49900.36050589249"
$ws.Range("R17").Value2 = [double]"449103.2445530324"
$ws.Range("S17").Value2 = [double]"0.2080055468919617"
$ws.Range("T17").Value2 = [double]"0.2080055468919617"
